# Append round-18 match rows (168-177 in col A; sheet rows 170-179) to leagueStats data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 10,15

$arr[0,0] = 168
$arr[0,1] = 'Cagliari'
$arr[0,2] = 'Inter'
$arr[0,3] = 0
$arr[0,4] = 3
$arr[0,5] = 0.31
$arr[0,6] = 3.52
$arr[0,7] = 0.36
$arr[0,8] = 3.56
$arr[0,9] = 0
$arr[0,10] = 2
$arr[0,11] = 0.05
$arr[0,12] = 0.04
$arr[0,13] = 0.09
$arr[0,14] = 1

$arr[1,0] = 169
$arr[1,1] = 'Empoli'
$arr[1,2] = 'Genoa'
$arr[1,3] = 1
$arr[1,4] = 2
$arr[1,5] = 1.34
$arr[1,6] = 1.77
$arr[1,7] = 1.22
$arr[1,8] = 1.77
$arr[1,9] = 1
$arr[1,10] = 1
$arr[1,11] = 0.12
$arr[1,12] = 0
$arr[1,13] = 0.13
$arr[1,14] = 1

$arr[2,0] = 170
$arr[2,1] = 'Lazio'
$arr[2,2] = 'Atalanta'
$arr[2,3] = 1
$arr[2,4] = 1
$arr[2,5] = 0.62
$arr[2,6] = 2.88
$arr[2,7] = 0.46
$arr[2,8] = 2.34
$arr[2,9] = 0
$arr[2,10] = 0
$arr[2,11] = 0.16
$arr[2,12] = 0.54
$arr[2,13] = 0.7
$arr[2,14] = 2

$arr[3,0] = 171
$arr[3,1] = 'Parma'
$arr[3,2] = 'Monza'
$arr[3,3] = 2
$arr[3,4] = 1
$arr[3,5] = 1.36
$arr[3,6] = 2.63
$arr[3,7] = 1.51
$arr[3,8] = 1.21
$arr[3,9] = 1
$arr[3,10] = 0
$arr[3,11] = 0.15
$arr[3,12] = 1.42
$arr[3,13] = 1.58
$arr[3,14] = 2

$arr[4,0] = 172
$arr[4,1] = 'Juventus'
$arr[4,2] = 'Fiorentina'
$arr[4,3] = 2
$arr[4,4] = 2
$arr[4,5] = 1.49
$arr[4,6] = 1
$arr[4,7] = 1.36
$arr[4,8] = 1.41
$arr[4,9] = 0
$arr[4,10] = 0
$arr[4,11] = 0.13
$arr[4,12] = 0.41
$arr[4,13] = 0.54
$arr[4,14] = 4

$arr[5,0] = 173
$arr[5,1] = 'Milan'
$arr[5,2] = 'Roma'
$arr[5,3] = 1
$arr[5,4] = 1
$arr[5,5] = 1.88
$arr[5,6] = 1.26
$arr[5,7] = 2.1
$arr[5,8] = 1.28
$arr[5,9] = 0
$arr[5,10] = 0
$arr[5,11] = 0.22
$arr[5,12] = 0.02
$arr[5,13] = 0.23
$arr[5,14] = 2

$arr[6,0] = 174
$arr[6,1] = 'Napoli'
$arr[6,2] = 'Venezia'
$arr[6,3] = 1
$arr[6,4] = 0
$arr[6,5] = 1.97
$arr[6,6] = 0.21
$arr[6,7] = 2.44
$arr[6,8] = 0.23
$arr[6,9] = 1
$arr[6,10] = 0
$arr[6,11] = 0.47
$arr[6,12] = 0.02
$arr[6,13] = 0.49
$arr[6,14] = 0

$arr[7,0] = 175
$arr[7,1] = 'Udinese'
$arr[7,2] = 'Torino'
$arr[7,3] = 2
$arr[7,4] = 2
$arr[7,5] = 0.68
$arr[7,6] = 0.6
$arr[7,7] = 0.79
$arr[7,8] = 0.78
$arr[7,9] = 0
$arr[7,10] = 0
$arr[7,11] = 0.11
$arr[7,12] = 0.18
$arr[7,13] = 0.29
$arr[7,14] = 4

$arr[8,0] = 176
$arr[8,1] = 'Bologna'
$arr[8,2] = 'Hellas Verona'
$arr[8,3] = 2
$arr[8,4] = 3
$arr[8,5] = 2.32
$arr[8,6] = 0.84
$arr[8,7] = 2.71
$arr[8,8] = 0.71
$arr[8,9] = 0
$arr[8,10] = 0
$arr[8,11] = 0.39
$arr[8,12] = 0.13
$arr[8,13] = 0.52
$arr[8,14] = 5

$arr[9,0] = 177
$arr[9,1] = 'Como'
$arr[9,2] = 'Lecce'
$arr[9,3] = 2
$arr[9,4] = 0
$arr[9,5] = 2.03
$arr[9,6] = 0.43
$arr[9,7] = 1.81
$arr[9,8] = 0.64
$arr[9,9] = 1
$arr[9,10] = 0
$arr[9,11] = 0.22
$arr[9,12] = 0.21
$arr[9,13] = 0.44
$arr[9,14] = 1

# Write the 10x15 block (rows 170-179, columns A-O) in one shot
$ws.Range("A170:O179").Value() = $arr

# Match the match-index column style (bold, thin border, centered) used for A2:A169
$ws.Range("A169").Copy()
$ws.Range("A170:A179").PasteSpecial(-4122)
$excel.CutCopyMode = 0

